$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text for the "test-19" row (B17)
$ws.Range("B17").Value = "RC, closed loop control && Servo"

# Move the active selection to B17
$ws.Range("B17").Select()
